$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> Efna1/Epha5 -> MuSCs) updated metrics ---
$ws.Range("G2").Value = 14.02618833333333
$ws.Range("H2").Value = 42.078565
$ws.Range("I2").Value = 0.806325281849088
$ws.Range("J2").Value = 0.8172785134657441
$ws.Range("M2").Value = 0.0237255
$ws.Range("N2").Value = 0.047451
$ws.Range("Q2").Value = 0.3327783313025
$ws.Range("R2").Value = 1.996669987815
$ws.Range("S2").Value = 0.806325281849088
$ws.Range("T2").Value = 0.8172785134657441

# --- Row 3 (FAPs -> Efna1/Epha5 -> MuSCs) updated metrics ---
$ws.Range("I3").Value = 0.1427547569137158
$ws.Range("J3").Value = 0.1446939568272663
$ws.Range("M3").Value = 0.0237255
$ws.Range("N3").Value = 0.047451
$ws.Range("Q3").Value = 0.058916284607
$ws.Range("R3").Value = 0.353497707642
$ws.Range("S3").Value = 0.1427547569137158
$ws.Range("T3").Value = 0.1446939568272663

# --- Row 4: sending cluster becomes the new "Inflammatory-Mac" cluster, with its own metrics ---
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1863673333333333
$ws.Range("H4").Value = 0.559102
$ws.Range("I4").Value = 0.01071372271683668
$ws.Range("J4").Value = 0.01085925937435662
$ws.Range("M4").Value = 0.0237255
$ws.Range("N4").Value = 0.047451
$ws.Range("Q4").Value = 0.004421658167
$ws.Range("R4").Value = 0.026529949002
$ws.Range("S4").Value = 0.01071372271683668
$ws.Range("T4").Value = 0.01085925937435662

# --- Row 5 (new): MuSCs -> Efna1/Epha5 -> MuSCs, freshly computed metrics ---
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Efna1"
$ws.Range("C5").Value = "Epha5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6993955000000001
$ws.Range("H5").Value = 1.398791
$ws.Range("I5").Value = 0.04020623852035952
$ws.Range("J5").Value = 0.02716827033263282
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0237255
$ws.Range("N5").Value = 0.047451
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.01659350793525
$ws.Range("R5").Value = 0.06637403174100001
$ws.Range("S5").Value = 0.04020623852035952
$ws.Range("T5").Value = 0.02716827033263282
